$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "HK_G_acc_LG"
$ws.Range("A2").Value = 76.51195499296765
$ws.Range("A3").Value = 76.51195499296765
$ws.Range("A4").Value = 76.23066104078762
$ws.Range("A5").Value = 76.51195499296765
$ws.Range("A6").Value = 76.37130801687763
$ws.Range("A7").Value = 76.37130801687763
$ws.Range("A8").Value = 76.37130801687763
$ws.Range("A9").Value = 76.09001406469761
$ws.Range("A10").Value = 76.37130801687763
$ws.Range("A11").Value = 76.37130801687763
$ws.Range("A12").Value = 75.9493670886076
$ws.Range("A13").Value = 76.51195499296765
$ws.Range("A14").Value = 76.51195499296765
$ws.Range("A15").Value = 76.23066104078762
$ws.Range("A16").Value = 76.51195499296765
$ws.Range("A17").Value = 76.51195499296765
$ws.Range("A18").Value = 76.51195499296765
$ws.Range("A19").Value = 76.51195499296765
$ws.Range("A20").Value = 76.23066104078762
$ws.Range("A21").Value = 76.23066104078762
$ws.Range("A22").Value = 76.23066104078762
$ws.Range("A23").Value = 76.09001406469761
$ws.Range("A24").Value = 76.23066104078762
$ws.Range("A25").Value = 75.80872011251758
$ws.Range("A26").Value = 77.35583684950772
$ws.Range("A27").Value = 76.51195499296765
$ws.Range("A28").Value = 77.35583684950772
$ws.Range("A29").Value = 76.09001406469761
$ws.Range("A30").Value = 76.23066104078762
$ws.Range("A31").Value = 76.51195499296765
$ws.Range("A32").Value = 76.37130801687763
$ws.Range("A33").Value = 76.51195499296765
$ws.Range("A34").Value = 76.51195499296765
$ws.Range("A35").Value = 76.65260196905767
$ws.Range("A36").Value = 75.66807313642757
$ws.Range("A37").Value = 75.38677918424754
$ws.Range("A38").Value = 76.09001406469761
$ws.Range("A39").Value = 77.0745428973277
$ws.Range("A40").Value = 77.21518987341773
$ws.Range("A41").Value = 76.51195499296765
$ws.Range("A42").Value = 76.37130801687763
$ws.Range("A43").Value = 76.37130801687763
$ws.Range("A44").Value = 76.37130801687763
$ws.Range("A45").Value = 76.51195499296765
$ws.Range("A46").Value = 76.37130801687763
$ws.Range("A47").Value = 76.51195499296765
$ws.Range("A48").Value = 76.37130801687763
$ws.Range("A49").Value = 76.51195499296765
